$d = $word.ActiveDocument

# 1) The paragraph "Cuando la cantidad de colores se vuelve mayor a 2, "
#    loses its text, leaving just a single space.
$null = $d.Content.Find.Execute(
    "Cuando la cantidad de colores se vuelve mayor a 2, ",
    $true, $false, $false, $false, $false,
    $true, 1, $false, " ", 2)

# 2) That paragraph is merged into the previous one ("Parte B:", which
#    carries the _GoBack bookmark) by removing the paragraph mark that
#    separates them. Locate the "Parte B:" paragraph dynamically so the
#    script does not depend on a hard-coded paragraph index.
$targetIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $paraText = $d.Paragraphs($i).Range.Text.TrimEnd([char]13)
    if ($paraText -eq "Parte B:") {
        $targetIndex = $i
        break
    }
}

if ($targetIndex -gt 0) {
    $partBPara = $d.Paragraphs($targetIndex)
    $markEnd = $partBPara.Range.End
    $markRange = $d.Range($markEnd - 1, $markEnd)
    $markRange.Delete()
}
